# MASSACHUSETTS_2023.xlsx cleanup
# 1. Rename header columns to short machine-friendly codes.
# 2. Convert the Estado/Municipio text columns (A & B) from ALL CAPS to
#    Proper/Title case for the data rows (2-393).
# 3. Remove the trailing footnote/metadata rows (395-399) that sat below
#    the data table, and shrink the used range back down to D393.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header row -------------------------------------------------------
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- 2. Title-case the state (A) and municipality (B) columns ------------
for ($r = 2; $r -le 393; $r++) {
    $aCell = $ws.Cells.Item($r, 1)
    $aVal = $aCell.Value2
    if ($aVal -ne $null -and $aVal -ne "") {
        $aCell.Value = $excel.WorksheetFunction.Proper($aVal)
    }

    $bCell = $ws.Cells.Item($r, 2)
    $bVal = $bCell.Value2
    if ($bVal -ne $null -and $bVal -ne "") {
        $bCell.Value = $excel.WorksheetFunction.Proper($bVal)
    }
}

# --- 3. Drop the trailing metadata/footnote rows (395-399) ---------------
$ws.Range("A395:D399").ClearContents()
